# Update countries & provincias Spain
#
# The "Pais" sheet lists countries with covid-style stats
# (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) sorted descending by
# column B (Casos totales). This refreshes the numbers for the
# rows whose figures moved and re-labels the rows whose relative
# rank swapped with a neighbour once the new totals are applied
# (Bielorrusia/Irlanda, Bolivia/Azerbaiyan, Suazilandia/Guadalupe/
# Togo, Nueva Caledonia/Belice, Islas Virgenes Britanicas/Butan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1352320
$ws.Range("C4").Value = 5011
$ws.Range("E4").Value = 1033149
$ws.Range("G4").Value = 286
$ws.Range("H4").Value = 80323

# Row 15 - Canada
$ws.Range("B15").Value = 68003
$ws.Range("C15").Value = 301
$ws.Range("D15").Value = 31644
$ws.Range("E15").Value = 31631
$ws.Range("G15").Value = 35
$ws.Range("H15").Value = 4728

# Row 29 - was Bielorrusia, now Irlanda (rank swap)
$ws.Range("A29").Value = "Irlanda"
$ws.Range("B29").Value = 22996
$ws.Range("C29").Value = 236
$ws.Range("D29").Value = 17110
$ws.Range("E29").Value = 4428
$ws.Range("F29").Value = 72
$ws.Range("G29").Value = 12
$ws.Range("H29").Value = 1458

# Row 30 - was Irlanda, now Bielorrusia (rank swap)
$ws.Range("A30").Value = "Bielorrusia"
$ws.Range("B30").Value = 22973
$ws.Range("C30").Value = 921
$ws.Range("D30").Value = 6406
$ws.Range("E30").Value = 16436
$ws.Range("F30").Value = 92
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 131

# Row 55 - Marruecos
$ws.Range("B55").Value = 6063
$ws.Range("C55").Value = 153
$ws.Range("D55").Value = 2554
$ws.Range("E55").Value = 3321

# Row 73 - was Bolivia, now Azerbaiyan (rank swap)
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 2519
$ws.Range("C73").Value = 97
$ws.Range("D73").Value = 1650
$ws.Range("E73").Value = 837
$ws.Range("F73").Value = 33
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 32

# Row 74 - was Azerbaiyan, now Bolivia (rank swap)
$ws.Range("A74").Value = "Bolivia"
$ws.Range("B74").Value = 2437
$ws.Range("C74").Value = 171
$ws.Range("D74").Value = 258
$ws.Range("E74").Value = 2065
$ws.Range("F74").Value = 3
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 114

# Row 121 - Jordania
$ws.Range("B121").Value = 540
$ws.Range("C121").Value = 18
$ws.Range("D121").Value = 389
$ws.Range("E121").Value = 142

# Row 127 - Reunion
$ws.Range("B127").Value = 436
$ws.Range("C127").Value = 5
$ws.Range("E127").Value = 82

# Row 148 - was Suazilandia, now Togo (3-way rank rotation)
$ws.Range("A148").Value = "Togo"
$ws.Range("B148").Value = 173
$ws.Range("C148").Value = 20
$ws.Range("D148").Value = 89
$ws.Range("E148").Value = 73
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 11

# Row 149 - was Guadalupe, now Suazilandia
$ws.Range("A149").Value = "Suazilandia"
$ws.Range("B149").Value = 163
$ws.Range("D149").Value = 14
$ws.Range("E149").Value = 147
$ws.Range("F149").Value = 0
$ws.Range("H149").Value = 2

# Row 150 - was Togo, now Guadalupe
$ws.Range("A150").Value = "Guadalupe"
$ws.Range("B150").Value = 154
$ws.Range("D150").Value = 104
$ws.Range("E150").Value = 37
$ws.Range("F150").Value = 4
$ws.Range("H150").Value = 13

# Row 192 - was Nueva Caledonia, now Belice (rank swap)
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Row 193 - was Belice, now Nueva Caledonia (rank swap)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

# Row 212 - was Islas Virgenes Britanicas, now Butan (rank swap)
$ws.Range("A212").Value = "Butan"
$ws.Range("D212").Value = 5
$ws.Range("H212").Value = 0

# Row 213 - was Butan, now Islas Virgenes Britanicas (rank swap)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 4
$ws.Range("H213").Value = 1
